$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

$ws.Range("H107").Value = 1439.5555
$ws.Range("I107").Value = 1432
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1432
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 488
$ws.Range("N107").Value = -5340

$ws.Range("H125").Value = 7721.375
$ws.Range("I125").Value = 7597.25
$ws.Range("J125").Value = 7845.5
$ws.Range("K125").Value = 68375.25
$ws.Range("L125").Value = 70609.5
$ws.Range("M125").Value = -65915.25
$ws.Range("N125").Value = -75529.5

$ws.Range("H132").Value = 2400.5
$ws.Range("I132").Value = 2400.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7201.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4671.5
$ws.Range("N132").ClearContents()

$ws.Range("H137").Value = 100
$ws.Range("I137").Value = 100
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 300
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 2250
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4598.353
$ws.Range("I32").Value = 4292.3438
$ws.Range("J32").Value = 9494.5
$ws.Range("K32").Value = 4292.3438
$ws.Range("L32").Value = 9494.5
$ws.Range("M32").Value = -4005.3438
$ws.Range("N32").Value = -10068.5

$ws.Range("H61").Value = 3121.4583
$ws.Range("I61").Value = 3321.875
$ws.Range("J61").Value = 2720.625
$ws.Range("K61").Value = 3321.875
$ws.Range("L61").Value = 2720.625
$ws.Range("M61").Value = -3109.875
$ws.Range("N61").Value = -3144.625

$ws.Range("H74").Value = 2860.7144
$ws.Range("I74").Value = 1150
$ws.Range("J74").Value = 3327.2727
$ws.Range("K74").Value = 1150
$ws.Range("L74").Value = 3327.2727
$ws.Range("M74").Value = -276

$ws.Range("H77").Value = 2860.7144
$ws.Range("I77").Value = 1150
$ws.Range("J77").Value = 3327.2727
$ws.Range("K77").Value = 5750
$ws.Range("L77").Value = 16636.3635
$ws.Range("M77").Value = -1382

$ws.Range("H132").Value = 2378.7742
$ws.Range("I132").Value = 2017.6818
$ws.Range("J132").Value = 3261.4443
$ws.Range("K132").Value = 6053.0454
$ws.Range("L132").Value = 9784.332900000001
$ws.Range("M132").Value = -3523.0454
$ws.Range("N132").Value = -14844.3329

$ws.Range("H136").Value = 3121.4583
$ws.Range("I136").Value = 3321.875
$ws.Range("J136").Value = 2720.625
$ws.Range("K136").Value = 9965.625
$ws.Range("L136").Value = 8161.875
$ws.Range("M136").Value = -7415.625
$ws.Range("N136").Value = -13261.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2752.96
$ws.Range("I20").Value = 2601.4285
$ws.Range("J20").Value = 2945.818
$ws.Range("K20").Value = 2601.4285
$ws.Range("L20").Value = 2945.818
$ws.Range("M20").Value = -2354.4285
$ws.Range("N20").Value = -3439.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5429.6523
$ws.Range("I31").Value = 2407.5454
$ws.Range("J31").Value = 8199.916999999999
$ws.Range("K31").Value = 2407.5454
$ws.Range("L31").Value = 8199.916999999999
$ws.Range("M31").Value = -2112.5454
$ws.Range("N31").Value = -8789.916999999999

$ws.Range("H34").Value = 5429.6523
$ws.Range("I34").Value = 2407.5454
$ws.Range("J34").Value = 8199.916999999999
$ws.Range("K34").Value = 2407.5454
$ws.Range("L34").Value = 8199.916999999999
$ws.Range("M34").Value = -2205.5454
$ws.Range("N34").Value = -8603.916999999999

$ws.Range("H58").Value = 976.3333
$ws.Range("I58").Value = 790.5
$ws.Range("J58").Value = 1348
$ws.Range("K58").Value = 790.5
$ws.Range("L58").Value = 1348
$ws.Range("M58").Value = -587.5
$ws.Range("N58").Value = -1754

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 3542.111
$ws.Range("I132").Value = 2751.2222
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 8253.6666
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -5723.6666
$ws.Range("N132").Value = -18059

$ws.Range("H134").Value = 1995.0588
$ws.Range("I134").Value = 1901.0667
$ws.Range("J134").Value = 2700
$ws.Range("K134").Value = 5703.2001
$ws.Range("L134").Value = 8100
$ws.Range("M134").Value = -3168.2001
$ws.Range("N134").Value = -13170

$ws.Range("H136").Value = 976.3333
$ws.Range("I136").Value = 790.5
$ws.Range("J136").Value = 1348
$ws.Range("K136").Value = 2371.5
$ws.Range("L136").Value = 4044
$ws.Range("M136").Value = 178.5
$ws.Range("N136").Value = -9144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4619.6
$ws.Range("I68").Value = 3399
$ws.Range("J68").Value = 4924.75
$ws.Range("K68").Value = 10197
$ws.Range("L68").Value = 14774.25
$ws.Range("M68").Value = -9386
$ws.Range("N68").Value = -16396.25

$ws.Range("H71").Value = 4619.6
$ws.Range("I71").Value = 3399
$ws.Range("J71").Value = 4924.75
$ws.Range("K71").Value = 30591
$ws.Range("L71").Value = 44322.75
$ws.Range("M71").Value = -26535
$ws.Range("N71").Value = -52434.75

$ws.Range("H122").Value = 3484.3877
$ws.Range("I122").Value = 480
$ws.Range("J122").Value = 3680.3262
$ws.Range("K122").Value = 4320
$ws.Range("L122").Value = 33122.9358
$ws.Range("M122").Value = -1870
$ws.Range("N122").Value = -38022.9358

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -888

$ws.Range("H132").Value = 2768.5264
$ws.Range("I132").Value = 2223.6155
$ws.Range("J132").Value = 3949.1667
$ws.Range("K132").Value = 6670.8465
$ws.Range("L132").Value = 11847.5001
$ws.Range("M132").Value = -4140.8465
$ws.Range("N132").Value = -16907.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 39998.332
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 39998.332
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 39998.332
$ws.Range("N54").Value = -41038.332

$ws.Range("H96").Value = 3570
$ws.Range("I96").Value = 3570
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 3570
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -2197
$ws.Range("N96").ClearContents()

$ws.Range("H101").Value = 37498.75
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 37498.75
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 37498.75
$ws.Range("N101").Value = -43988.75

$ws.Range("H132").Value = 2502.8667
$ws.Range("I132").Value = 2144.0293
$ws.Range("J132").Value = 3612
$ws.Range("K132").Value = 6432.0879
$ws.Range("L132").Value = 10836
$ws.Range("M132").Value = -3902.0879
$ws.Range("N132").Value = -15896

$ws.Range("H136").Value = 11853.833
$ws.Range("I136").Value = 17298.715
$ws.Range("J136").Value = 4231
$ws.Range("K136").Value = 51896.145
$ws.Range("L136").Value = 12693
$ws.Range("M136").Value = -49346.145
$ws.Range("N136").Value = -17793
